$d = $word.ActiveDocument

# Hybrid bold + color (#2C3E50) highlighting for quantitative impact metrics.
# Word's Font.Color uses a BGR-packed long (wdColor), so RGB 2C3E50 ->
# 0x503E2C = 5258796.
$metricColor = 5258796

function Apply-MetricHighlights($Paragraph, $AnchorText, $Tokens) {
    $full = $Paragraph.Range.Text
    if ($full.IndexOf($AnchorText) -lt 0) {
        Write-Warning ("Paragraph text did not match expected anchor: " + $AnchorText)
    }

    $pEnd = $Paragraph.Range.End
    $searchStart = $Paragraph.Range.Start

    foreach ($tok in $Tokens) {
        $r = $d.Range($searchStart, $pEnd)
        $found = $r.Find.Execute($tok, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $r.Font.Bold = 1
            $r.Font.Color = $metricColor
            $searchStart = $r.End
        }
    }
}

# "• Discovered systematic race coding errors ... from 23% to 64%"
Apply-MetricHighlights $d.Paragraphs(9) "Discovered systematic race coding errors" @("23%", "64%")

# "• Achieved 87% prediction accuracy ... of 71%, reducing polling error margins from ±4.2% to ±2.1%"
Apply-MetricHighlights $d.Paragraphs(11) "reducing polling error margins" @("87%", "71%", "±4.2%", "±2.1%")

# "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
Apply-MetricHighlights $d.Paragraphs(31) "Wrote RFP and analyzed bids" @("1,200")

# "• Created comprehensive meta-analysis framework ... became the $400M Polling Consortium Database ... now valued at $1B+"
Apply-MetricHighlights $d.Paragraphs(46) "Polling Consortium Database" @("$400M", "$1B")

# "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
Apply-MetricHighlights $d.Paragraphs(63) "Algorithm reduced mapping costs" @("73.5%", "$4.7M")

# "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"
Apply-MetricHighlights $d.Paragraphs(65) "industry standard of" @("87%", "71%")

Write-Output "Applied metric highlights"
